# Populate the pharmacy "transactions" report with the full data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 13 additional data rows after row 4 -------------
# (row 5 currently holds the blank subtotal placeholder, row 6 the footer;
#  both need to end up 13 rows further down, at rows 18 and 19.)
$ws.Rows("5:17").Insert()

# --- 2. Clone row 4's cell formatting onto the freshly inserted rows ------
$ws.Range("A4:N4").Copy()
$ws.Range("A5:N17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Re-create the merges that PasteSpecial(Formats) does not copy ----
for ($r = 5; $r -le 17; $r++) {
    $ws.Range("B$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
}

# --- 4. Fill in the data rows ----------------------------------------------
$data = @(
    @{ n = 1;  name = "ANGIOFOX (EFFOX) 25MG LONG 30 CAPS.";       ratio = "0:0";    qty = 114;    pct = 1 },
    @{ n = 2;  name = "AUGMENTIN 457MG/5ML SUSP. 70 ML";           ratio = "1:0";    qty = 137;    pct = 1 },
    @{ n = 3;  name = "BLOKATENS 10/160MG 28 F.C.TABS.";           ratio = "0:0";    qty = 160;    pct = 1 },
    @{ n = 4;  name = "COLOVATIL 30 F.C. TABS";                    ratio = "0:0";    qty = 63;     pct = 1 },
    @{ n = 5;  name = "GAVISCON LIQUID 24 SACHETS 10 ML";          ratio = "0:9";    qty = 12;     pct = 0.04 },
    @{ n = 6;  name = "GINKGO BILOBA 30 CAPS.";                    ratio = "0:0";    qty = 186;    pct = 1 },
    @{ n = 7;  name = "MILGA ADVANCE 30 F.C. TABS";                ratio = "0:0";    qty = 136.5;  pct = 1 },
    @{ n = 8;  name = "PERLOC 40MG 14 F.C.TAB.";                   ratio = "0:0";    qty = 68.25;  pct = 1 },
    @{ n = 9;  name = "RHINEX 0.05% INFANTILE NASAL DROPS 10 ML";  ratio = "2:0";    qty = 18;     pct = 1 },
    @{ n = 10; name = "RIVO 320MG 20*10 TABS";                     ratio = "1:2";    qty = 14.1;   pct = 0.1 },
    @{ n = 11; name = "VASTAREL MR 35MG 30 F.C.TAB.";              ratio = "2:0";    qty = 175;    pct = 1 },
    @{ n = 12; name = "WATER FOR INJECTION AMP. 5 ML";             ratio = "7816:0"; qty = 2.5;    pct = 1 },
    @{ n = 13; name = "سويت كوكو";                                 ratio = "22:0";   qty = 25;     pct = 1 },
    @{ n = 14; name = "مرطب شفاه لونا جوز هند ابيض";                ratio = "3:0";    qty = 20;     pct = 1 }
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.n
    $ws.Range("B$row").Value = $item.name
    $ws.Range("H$row").Value = $item.ratio
    $ws.Range("L$row").Value = $item.qty
    $ws.Range("N$row").Value = $item.pct
    $row = $row + 1
}

# --- 5. Totals row (was the blank "K5:N5" row, now shifted to row 18) -----
$ws.Range("K18").Value = 1131.3499999999999

# --- 6. Footer row (was row 6, shifted to row 19) keeps its original values,
#        which the row-insert already carried down automatically.

Write-Output "populated transactions report"
